$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "37.028.43"
Set-TextValue $ws.Range("E2") "  -0.64%  "
Set-TextValue $ws.Range("D3") "2.007.78"
Set-TextValue $ws.Range("E3") "  -1.61%  "
Set-TextValue $ws.Range("E4") "  -0.13%  "
Set-TextValue $ws.Range("D5") "225.41"
Set-TextValue $ws.Range("E5") "  -1.28%  "
Set-TextValue $ws.Range("E6") "  -0.67%  "
Set-TextValue $ws.Range("E7") "  +0.01%  "
Set-TextValue $ws.Range("D8") "55.03"
Set-TextValue $ws.Range("E8") "  -1.74%  "
Set-TextValue $ws.Range("D9") "0.373"
Set-TextValue $ws.Range("E9") "  -2.77%  "
Set-TextValue $ws.Range("E10") "  -4.11%  "
Set-TextValue $ws.Range("E11") "  -4.50%  "
Set-TextValue $ws.Range("D12") "2.304.62"
Set-TextValue $ws.Range("E12") "  -1.50%  "
Set-TextValue $ws.Range("D13") "13.96"
Set-TextValue $ws.Range("E13") "  -3.80%  "
Set-TextValue $ws.Range("D14") "19.65"
Set-TextValue $ws.Range("E14") "  -4.12%  "
Set-TextValue $ws.Range("E15") "  -1.96%  "
Set-TextValue $ws.Range("D16") "0.732"
Set-TextValue $ws.Range("E16") "  -2.55%  "
Set-TextValue $ws.Range("D17") "1.997.56"
Set-TextValue $ws.Range("E17") "  -1.84%  "
Set-TextValue $ws.Range("D18") "36.978.76"
Set-TextValue $ws.Range("E18") "  -0.35%  "
Set-TextValue $ws.Range("D19") "6.14"
Set-TextValue $ws.Range("E19") "  +3.09%  "
Set-TextValue $ws.Range("D20") "68.14"
Set-TextValue $ws.Range("E20") "  -2.00%  "
Set-TextValue $ws.Range("D21") "0.0₃0807"
Set-TextValue $ws.Range("E21") "  -3.74%  "
Set-TextValue $ws.Range("D22") "223.40"
Set-TextValue $ws.Range("E22") "  -0.85%  "
Set-TextValue $ws.Range("E23") "  -0.01%  "
Set-TextValue $ws.Range("E24") "  +2.37%  "
Set-TextValue $ws.Range("D25") "2.15"
Set-TextValue $ws.Range("E25") "  -5.36%  "
Set-TextValue $ws.Range("D26") "163.98"
Set-TextValue $ws.Range("E26") "  -2.38%  "
Set-TextValue $ws.Range("D27") "8.88"
Set-TextValue $ws.Range("E27") "  -6.38%  "
Set-TextValue $ws.Range("E28") "  -1.85%  "
Set-TextValue $ws.Range("E29") "  -3.88%  "
Set-TextValue $ws.Range("D30") "1.29"
Set-TextValue $ws.Range("E30") "  -6.62%  "
Set-TextValue $ws.Range("E31") "  -1.44%  "
Set-TextValue $ws.Range("D32") "4.38"
Set-TextValue $ws.Range("E32") "  -2.48%  "
Set-TextValue $ws.Range("E33") "  -1.96%  "
Set-TextValue $ws.Range("D34") "4.45"
Set-TextValue $ws.Range("E34") "  -1.91%  "
Set-TextValue $ws.Range("B35") "LidoDAOToken"
Set-TextValue $ws.Range("C35") "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws.Range("D35") "2.30"
Set-TextValue $ws.Range("E35") "  -3.54%  "
Set-TextValue $ws.Range("B36") "WEMIXToken"
Set-TextValue $ws.Range("C36") "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D36") "1.86"
Set-TextValue $ws.Range("E36") "  +2.49%  "
Set-TextValue $ws.Range("E37") "  -0.10%  "
Set-TextValue $ws.Range("E38") "  -2.58%  "
Set-TextValue $ws.Range("D39") "5.33"
Set-TextValue $ws.Range("E39") "  -0.75%  "
Set-TextValue $ws.Range("D40") "1.456.52"
Set-TextValue $ws.Range("E40") "  -2.15%  "
Set-TextValue $ws.Range("B41") "VeChain"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D41") "0.0211"
Set-TextValue $ws.Range("E41") "  -3.71%  "
Set-TextValue $ws.Range("B42") "Aave"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D42") "94.37"
Set-TextValue $ws.Range("E42") "  -0.89%  "
Set-TextValue $ws.Range("B43") "FTXToken"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue $ws.Range("D43") "4.26"
Set-TextValue $ws.Range("E43") "  +17.82%  "
Set-TextValue $ws.Range("B44") "Cronos"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D44") "0.0906"
Set-TextValue $ws.Range("E44") "  -2.78%  "
Set-TextValue $ws.Range("B45") "HuobiToken"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D45") "2.75"
Set-TextValue $ws.Range("E45") "  -4.43%  "
Set-TextValue $ws.Range("D46") "15.87"
Set-TextValue $ws.Range("E46") "  -5.11%  "
Set-TextValue $ws.Range("E47") "  -2.75%  "
Set-TextValue $ws.Range("D48") "0.994"
Set-TextValue $ws.Range("E48") "  -1.57%  "
Set-TextValue $ws.Range("D49") "7.04"
Set-TextValue $ws.Range("E49") "  -0.67%  "
Set-TextValue $ws.Range("E50") "  -0.58%  "
Set-TextValue $ws.Range("D51") "2.192.97"
Set-TextValue $ws.Range("E51") "  -1.62%  "
